$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45178 = 2023-09-09) for
# every data row (rows 2-150). The update bumps that date forward by one day
# to serial 45179 (2023-09-10) for all of them.
for ($r = 2; $r -le 150; $r++) {
    $ws.Cells.Item($r, 3).Value = 45179
}
